# Force a full recalculation so the volatile RANDBETWEEN()-based RES
# generation formulas (Pg, Winter/Summer S1-S3) redraw new random values,
# then restore the active sheet/selection to mirror the author's click on
# "RES installed"!G22.

$wb = $excel.ActiveWorkbook
$excel.CalculateFullRebuild()

$ws = $wb.Worksheets.Item("RES installed")
$ws.Activate()
$ws.Range("G22").Select()
